$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5847.2856
$ws.Range("I74").Value = 6476.2
$ws.Range("J74").Value = 4275
$ws.Range("K74").Value = 6476.2
$ws.Range("L74").Value = 4275
$ws.Range("M74").Value = -5540.2
$ws.Range("N74").Value = -6147

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 5847.2856
$ws.Range("I77").Value = 6476.2
$ws.Range("J77").Value = 4275
$ws.Range("K77").Value = 32381
$ws.Range("L77").Value = 21375
$ws.Range("M77").Value = -27701
$ws.Range("N77").Value = -30735

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4640.1
$ws.Range("I86").Value = 6480.2
$ws.Range("J86").Value = 2800
$ws.Range("K86").Value = 6480.2
$ws.Range("L86").Value = 2800
$ws.Range("M86").Value = -5357.2
$ws.Range("N86").Value = -5046

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4640.1
$ws.Range("I89").Value = 6480.2
$ws.Range("J89").Value = 2800
$ws.Range("K89").Value = 32401
$ws.Range("L89").Value = 14000
$ws.Range("M89").Value = -26785
$ws.Range("N89").Value = -25232

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1013664.6
$ws.Range("I113").Value = 2780027.8
$ws.Range("J113").Value = 4314.2856
$ws.Range("K113").Value = 2780027.8
$ws.Range("L113").Value = 4314.2856
$ws.Range("M113").Value = -2776773.8
$ws.Range("N113").Value = -10822.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2005751.8
$ws.Range("I116").Value = 7145713.5
$ws.Range("J116").Value = 6877.778
$ws.Range("K116").Value = 7145713.5
$ws.Range("L116").Value = 6877.778
$ws.Range("M116").Value = -7142271.5
$ws.Range("N116").Value = -13761.778

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2978081.2
$ws.Range("I132").Value = 3827617.8
$ws.Range("K132").Value = 11482853.4
$ws.Range("M132").Value = -11480323.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4272.9253
$ws.Range("I138").Value = 1551.1364
$ws.Range("J138").Value = 5603.5776
$ws.Range("K138").Value = 4653.4092
$ws.Range("L138").Value = 16810.7328
$ws.Range("M138").Value = 486.5907999999999
$ws.Range("N138").Value = -27090.7328

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2687.4048
$ws.Range("I61").Value = 1861.3462
$ws.Range("J61").Value = 4029.75
$ws.Range("K61").Value = 1861.3462
$ws.Range("L61").Value = 4029.75
$ws.Range("M61").Value = -1649.3462
$ws.Range("N61").Value = -4453.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1724.4828
$ws.Range("I102").Value = 1667.037
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1667.037
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -45.03700000000003
$ws.Range("N102").Value = -5744

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5001802
$ws.Range("I122").Value = 6251623
$ws.Range("J122").Value = 2520
$ws.Range("K122").Value = 18754869
$ws.Range("L122").Value = 7560
$ws.Range("M122").Value = -18752419
$ws.Range("N122").Value = -12460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2687.4048
$ws.Range("I136").Value = 1861.3462
$ws.Range("J136").Value = 4029.75
$ws.Range("K136").Value = 5584.0386
$ws.Range("L136").Value = 12089.25
$ws.Range("M136").Value = -3034.0386
$ws.Range("N136").Value = -17189.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2489.8096
$ws.Range("I86").Value = 2023.1538
$ws.Range("K86").Value = 2023.1538
$ws.Range("M86").Value = -900.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2489.8096
$ws.Range("I89").Value = 2023.1538
$ws.Range("K89").Value = 10115.769
$ws.Range("M89").Value = -4499.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5210108.5
$ws.Range("I134").Value = 6251459
$ws.Range("J134").Value = 3356.25
$ws.Range("K134").Value = 18754377
$ws.Range("L134").Value = 10068.75
$ws.Range("M134").Value = -18751842
$ws.Range("N134").Value = -15138.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2718.9443
$ws.Range("I62").Value = 2931
$ws.Range("J62").Value = 2637.3845
$ws.Range("K62").Value = 2931
$ws.Range("L62").Value = 2637.3845
$ws.Range("M62").Value = -2307
$ws.Range("N62").Value = -3885.3845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2718.9443
$ws.Range("I65").Value = 2931
$ws.Range("J65").Value = 2637.3845
$ws.Range("K65").Value = 14655
$ws.Range("L65").Value = 13186.9225
$ws.Range("M65").Value = -11535
$ws.Range("N65").Value = -19426.9225

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3025
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 3933.3333
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 11799.9999
$ws.Range("M22").Value = -731
$ws.Range("N22").Value = -12137.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 3025
$ws.Range("I27").Value = 300
$ws.Range("J27").Value = 3933.3333
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 11799.9999
$ws.Range("M27").Value = -798
$ws.Range("N27").Value = -12003.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 5661.6665
$ws.Range("I112").Value = 2000
$ws.Range("J112").Value = 5994.5454
$ws.Range("K112").Value = 6000
$ws.Range("L112").Value = 17983.6362
$ws.Range("M112").Value = -4892
$ws.Range("N112").Value = -20199.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 499.13333
$ws.Range("I113").Value = 422
$ws.Range("J113").Value = 527.1818
$ws.Range("K113").Value = 1266
$ws.Range("L113").Value = 1581.5454
$ws.Range("M113").Value = 904
$ws.Range("N113").Value = -5921.5454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 4166.6665
$ws.Range("I129").Value = 500
$ws.Range("K129").Value = 1500
$ws.Range("M129").Value = 3500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 19308760
$ws.Range("I131").Value = 55555844
$ws.Range("J131").Value = 119127.12
$ws.Range("K131").Value = 166667532
$ws.Range("L131").Value = 357381.36
$ws.Range("M131").Value = -166662492
$ws.Range("N131").Value = -367461.36

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5333.3335
$ws.Range("I70").Value = 5333.3335
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 5333.3335
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -5063.3335
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5333.3335
$ws.Range("I73").Value = 5333.3335
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 5333.3335
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -4397.3335
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2868.923
$ws.Range("I122").Value = 2945.111
$ws.Range("J122").Value = 2697.5
$ws.Range("K122").Value = 8835.332999999999
$ws.Range("L122").Value = 8092.5
$ws.Range("M122").Value = -6385.332999999999
$ws.Range("N122").Value = -12992.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6541010
$ws.Range("I132").Value = 9013391
$ws.Range("J132").Value = 6860
$ws.Range("K132").Value = 27040173
$ws.Range("L132").Value = 20580
$ws.Range("M132").Value = -27037643
$ws.Range("N132").Value = -25640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 9898.5
$ws.Range("J43").Value = 9200
$ws.Range("L43").Value = 9200
$ws.Range("N43").Value = -9586

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5538.4116
$ws.Range("I122").Value = 5677.613
$ws.Range("J122").Value = 4100
$ws.Range("K122").Value = 17032.839
$ws.Range("L122").Value = 12300
$ws.Range("M122").Value = -17200
